# Bond dates update: "today" reference date moved forward by one day
# (2023-09-24 -> 2023-09-25). This shifts the cached "days since previous
# payout" (column G) up by 1 and the cached "days until next payout"
# (column I) down by 1, for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 262

for ($row = 2; $row -le $lastRow; $row++) {
    $gCell = $ws.Cells.Item($row, 7)   # column G: Dni od poprzedniej wypłaty
    $iCell = $ws.Cells.Item($row, 9)   # column I: Dni do następnej wypłaty

    $gVal = $gCell.Value()
    if ($gVal -ne $null) {
        $gCell.Value = $gVal + 1
    }

    $iVal = $iCell.Value()
    if ($iVal -ne $null) {
        $iCell.Value = $iVal - 1
    }
}
